# Wasstraat_Config_Harmonize.xlsx -- "Merge Artefacts toegevoegd: Cycle complete"
#
# 1. Spoor: add spooraard / beschrijving / datering rows
# 2. Vondst: add doosnr row
# 3. New sheet "Glas" (after Vondst) with glassoort / kleur / decoratie attributes
# 4. Artefact: give a few of the existing column-mapping lists an extra alias
#    (merge artifacts: "4b","7a","7b","9","10b" each gain a second, human
#    readable alias used by the merged table)
# 5. New sheet "Hout" (after Glas) with maten / houtsoortcd / bewerkingssporen /
#    gebruikssporen / decoratie attributes
# 6. Fix up final tab/selection state: Glas ends up the active tab, Artefact's
#    old D9 selection moves to B15, Vondst's selection collapses to A1:B1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Spoor: 3 new attribute rows
# ---------------------------------------------------------------------------
$spoor = $wb.Worksheets.Item("Spoor")
$spoor.Range("A5").Value = "spooraard"
$spoor.Range("B5").Value = '["AARD"]'
$spoor.Range("A6").Value = "beschrijving"
$spoor.Range("B6").Value = '["BESCHRIJF"]'
$spoor.Range("A7").Value = "datering"
$spoor.Range("B7").Value = '["DAT", "DATERING"]'
$spoor.Select()
$spoor.Range("A6:B6").Select()

# ---------------------------------------------------------------------------
# 2. Vondst: 1 new attribute row
# ---------------------------------------------------------------------------
$vondst = $wb.Worksheets.Item("Vondst")
$vondst.Range("A15").Value = "doosnr"
$vondst.Range("B15").Value = '["DOOSNR"]'

# ---------------------------------------------------------------------------
# 3. New sheet "Glas", inserted right after Vondst
# ---------------------------------------------------------------------------
$glas = $wb.Worksheets.Add($null, $vondst)
$glas.Name = "Glas"
$glas.Range("A1").Value = "Attribute"
$glas.Range("B1").Value = "Kolommen"
$glas.Range("A2").Value = "glassoort"
$glas.Range("B2").Value = "TEMP_GLASSOORT"
$glas.Range("A3").Value = "kleur"
$glas.Range("B3").Value = '["kleur"]'
$glas.Range("A4").Value = "decoratie"
$glas.Range("B4").Value = '["decoratie"]'
$glas.Columns.Item(2).ColumnWidth = 22.498697916666668

# ---------------------------------------------------------------------------
# 4. Artefact: extend a handful of column-mapping lists with a second alias
# ---------------------------------------------------------------------------
$artefact = $wb.Worksheets.Item("Artefact")
$artefact.Range("B9").Value = '["4b", "BESCHR"]'
$artefact.Range("B16").Value = '["7b", "FUNCTIE"]'
$artefact.Range("B18").Value = '["9", "LITERATUUR"]'
$artefact.Range("B19").Value = '["10b", "TEKNO"]'
$artefact.Range("B15").Value = '["7a", "TYPE"]'

# ---------------------------------------------------------------------------
# 5. New sheet "Hout", inserted right after Glas
# ---------------------------------------------------------------------------
$hout = $wb.Worksheets.Add($null, $glas)
$hout.Name = "Hout"
$hout.Range("A1").Value = "Attribute"
$hout.Range("B1").Value = "Kolommen"
$hout.Range("A2").Value = "maten"
$hout.Range("B2").Value = '["4a"]'
$hout.Range("A3").Value = "houtsoortcd"
$hout.Range("B3").Value = '["5a"]'
$hout.Range("A4").Value = "bewerkingssporen"
$hout.Range("B4").Value = '["5b"]'
$hout.Range("A5").Value = "gebruikssporen"
$hout.Range("B5").Value = '["5b1"]'
$hout.Range("A6").Value = "decoratie"
$hout.Range("B6").Value = '["5c"]'
$hout.Columns.Item(1).ColumnWidth = 15.166666666666666
$hout.Select()
$hout.Range("A7").Select()

# Go back and give Glas!B2 its real, final value (this is the very last new
# shared string created, matching the authoring order captured in the diff)
$glas.Range("B2").Value = '["glassoort", "GLSSOORT"]'

# ---------------------------------------------------------------------------
# 6. Final tab/selection clean-up
# ---------------------------------------------------------------------------
$artefact.Select()
$artefact.Range("B15").Select()

$vondst.Select()
$vondst.Range("A1:B1").Select()

$glas.Select()
$glas.Range("B2").Select()
